$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows to delete (original 1-indexed row numbers), senators without a
# recorded Twitter handle. Delete from the bottom up so row numbers of
# rows still to be removed don't shift.
$rowsToDelete = @(47, 45, 44, 43, 41, 36, 30, 24, 23, 13, 6, 3)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

$ws.Range("A36").Select()
